$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3196677091339666
$ws.Range("C2").Value = 0.04339544069544843
$ws.Range("E2").Value = 0.6731681490162913
$ws.Range("F2").Value = 2.276559655965301
$ws.Range("G2").Value = 0.002430289688350901
$ws.Range("I2").Value = 0.5037492469107399
$ws.Range("J2").Value = 0.04543404102100546
$ws.Range("K2").Value = 0.3286257949678486
$ws.Range("N2").Value = 1.31171027592978
$ws.Range("O2").Value = 2.330852130569127

$ws.Range("B3").Value = 0.2843165961010925
$ws.Range("C3").Value = 0.03788059928024268
$ws.Range("E3").Value = 0.645254848475588
$ws.Range("F3").Value = 2.24679625162679
$ws.Range("G3").Value = 0.002432391713797773
$ws.Range("I3").Value = 0.5095548468417874
$ws.Range("J3").Value = 0.04582909655639433
$ws.Range("K3").Value = 0.2900573016790986
$ws.Range("N3").Value = 1.324549781672978
$ws.Range("O3").Value = 2.350141138592591

$ws.Range("B4").Value = 0.2626137664573207
$ws.Range("C4").Value = 0.03448174982037244
$ws.Range("E4").Value = 0.628453221452574
$ws.Range("F4").Value = 2.229863603246784
$ws.Range("G4").Value = 0.002433750412389983
$ws.Range("I4").Value = 0.5133955085443631
$ws.Range("J4").Value = 0.04609840770437401
$ws.Range("K4").Value = 0.2663568733037494
$ws.Range("N4").Value = 1.332895341091536
$ws.Range("O4").Value = 2.363349020189617

$ws.Range("B5").Value = 0.2537709512309334
$ws.Range("C5").Value = 0.03309353959436123
$ws.Range("E5").Value = 0.6216911280433237
$ws.Range("F5").Value = 2.223300561978917
$ws.Range("G5").Value = 0.002434321254630167
$ws.Range("I5").Value = 0.5150299698095271
$ws.Range("J5").Value = 0.04621487299767857
$ws.Range("K5").Value = 0.2566944496635841
$ws.Range("N5").Value = 1.336412455406798
$ws.Range("O5").Value = 2.369074177446691

$ws.Range("B6").Value = 0.2523027006300822
$ws.Range("C6").Value = 0.03286283935403844
$ws.Range("E6").Value = 0.6205734008105566
$ws.Range("F6").Value = 2.222231126725831
$ws.Range("G6").Value = 0.002434417080644394
$ws.Range("I6").Value = 0.5153055588335071
$ws.Range("J6").Value = 0.04623461758017555
$ws.Range("K6").Value = 0.2550897669673589
$ws.Range("N6").Value = 1.337003490993236
$ws.Range("O6").Value = 2.370045535067518

$ws.Range("B7").Value = 0.2624945033302311
$ws.Range("C7").Value = 0.0344630406309534
$ws.Range("E7").Value = 0.6283616826546137
$ws.Range("F7").Value = 2.229773727145229
$ws.Range("G7").Value = 0.002433758041527045
$ws.Range("I7").Value = 0.5134172706790245
$ws.Range("J7").Value = 0.04609995120023136
$ws.Range("K7").Value = 0.266226579219591
$ws.Range("N7").Value = 1.332942303450533
$ws.Range("O7").Value = 2.363424843844484

$ws.Range("B8").Value = 0.3074783544509501
$ws.Range("C8").Value = 0.04149658802610645
$ws.Range("E8").Value = 0.66347357630913
$ws.Range("F8").Value = 2.266018428670563
$ws.Range("G8").Value = 0.002431000374511178
$ws.Range("I8").Value = 0.5056937090266835
$ws.Range("J8").Value = 0.04556470118813749
$ws.Range("K8").Value = 0.3153316705864881
$ws.Range("N8").Value = 1.316041466109088
$ws.Range("O8").Value = 2.337219663813102

$ws.Range("B9").Value = 0.3956952997791348
$ws.Range("C9").Value = 0.05518711476585736
$ws.Range("E9").Value = 0.7350126312320668
$ws.Range("F9").Value = 2.34776701983256
$ws.Range("G9").Value = 0.002426130185623543
$ws.Range("I9").Value = 0.4927395964035739
$ws.Range("J9").Value = 0.04472760018300193
$ws.Range("K9").Value = 0.4114556756525189
$ws.Range("N9").Value = 1.286562868561028
$ws.Range("O9").Value = 2.296668556491923

$ws.Range("B10").Value = 0.4604910639335458
$ws.Range("C10").Value = 0.06518236297753788
$ws.Range("E10").Value = 0.7892274922030964
$ws.Range("F10").Value = 2.414376137777765
$ws.Range("G10").Value = 0.002422876512005626
$ws.Range("I10").Value = 0.4845608227211358
$ws.Range("J10").Value = 0.04424259053751101
$ws.Range("K10").Value = 0.4819558461101678
$ws.Range("N10").Value = 1.267135028194719
$ws.Range("O10").Value = 2.273499638586898

$ws.Range("B11").Value = 0.4899608423780251
$ws.Range("C11").Value = 0.06971562052022762
$ws.Range("E11").Value = 0.814254998374409
$ws.Range("F11").Value = 2.446110000107126
$ws.Range("G11").Value = 0.002421466086915623
$ws.Range("I11").Value = 0.4811312906421463
$ws.Range("J11").Value = 0.04405027664324734
$ws.Range("K11").Value = 0.5139982935104968
$ws.Range("N11").Value = 1.258780333982834
$ws.Range("O11").Value = 2.264401704047742

$ws.Range("B12").Value = 0.501118937803966
$ws.Range("C12").Value = 0.0714302546095098
$ws.Range("E12").Value = 0.8237849555028873
$ws.Range("F12").Value = 2.458333446654365
$ws.Range("G12").Value = 0.002420941964008004
$ws.Range("I12").Value = 0.4798745173049497
$ws.Range("J12").Value = 0.04398153284152073
$ws.Range("K12").Value = 0.5261274038255976
$ws.Range("N12").Value = 1.255686071788105
$ws.Range("O12").Value = 2.26116418299182

$ws.Range("B13").Value = 0.4987159170351845
$ws.Range("C13").Value = 0.07106106779804122
$ws.Range("E13").Value = 0.8217301682369822
$ws.Range("F13").Value = 2.455691716388401
$ws.Range("G13").Value = 0.002421054400445348
$ws.Range("I13").Value = 0.480143320652072
$ws.Range("J13").Value = 0.04399615638838128
$ws.Range("K13").Value = 0.52351539979432
$ws.Range("N13").Value = 1.256349386282089
$ws.Range("O13").Value = 2.261852199322902

$ws.Range("B14").Value = 0.4908788584272941
$ws.Range("C14").Value = 0.069856725122321
$ws.Range("E14").Value = 0.815037979176239
$ws.Range("F14").Value = 2.447111488807451
$ws.Range("G14").Value = 0.002421422767225234
$ws.Range("I14").Value = 0.4810270547057129
$ws.Range("J14").Value = 0.04404453920261986
$ws.Range("K14").Value = 0.514996259613298
$ws.Range("N14").Value = 1.258524374658492
$ws.Range("O14").Value = 2.264131187223455

$ws.Range("B15").Value = 0.4860782248745465
$ws.Range("C15").Value = 0.06911876687645702
$ws.Range("E15").Value = 0.8109456714637986
$ws.Range("F15").Value = 2.441882760357402
$ws.Range("G15").Value = 0.002421649700431162
$ws.Range("I15").Value = 0.4815738279931523
$ws.Range("J15").Value = 0.04407470684577675
$ws.Range("K15").Value = 0.5097774149999736
$ws.Range("N15").Value = 1.259865666310304
$ws.Range("O15").Value = 2.265554188236592

$ws.Range("B16").Value = 0.4585649708464814
$ws.Range("C16").Value = 0.06488582519327224
$ws.Range("E16").Value = 0.7875992379799612
$ws.Range("F16").Value = 2.412331130241938
$ws.Range("G16").Value = 0.002422970084942784
$ws.Range("I16").Value = 0.484790812659142
$ws.Range("J16").Value = 0.04425572933015332
$ws.Range("K16").Value = 0.479861175028276
$ws.Range("N16").Value = 1.267690749319385
$ws.Range("O16").Value = 2.274123251158343

$ws.Range("B17").Value = 0.4416844900978276
$ws.Range("C17").Value = 0.06228552426878764
$ws.Range("E17").Value = 0.7733704947080753
$ws.Range("F17").Value = 2.394569518881781
$ws.Range("G17").Value = 0.002423797913512292
$ws.Range("I17").Value = 0.4868389158694235
$ws.Range("J17").Value = 0.04437404043146032
$ws.Range("K17").Value = 0.4615008350597805
$ws.Range("N17").Value = 1.272614928215127
$ws.Range("O17").Value = 2.279749603085932

$ws.Range("B18").Value = 0.4319747412033621
$ws.Range("C18").Value = 0.06078862003779761
$ws.Range("E18").Value = 0.7652208366309736
$ws.Range("F18").Value = 2.384488367072009
$ws.Range("G18").Value = 0.00242428062069059
$ws.Range("I18").Value = 0.4880443201234783
$ws.Range("J18").Value = 0.04444475486579513
$ws.Range("K18").Value = 0.4509378114484548
$ws.Range("N18").Value = 1.27549265939383
$ws.Range("O18").Value = 2.283121412548908

$ws.Range("B19").Value = 0.4286871117356839
$ws.Range("C19").Value = 0.06028157574826309
$ws.Range("E19").Value = 0.7624673967507647
$ws.Range("F19").Value = 2.381098210448727
$ws.Range("G19").Value = 0.002424445185609116
$ws.Range("I19").Value = 0.488457151457748
$ws.Range("J19").Value = 0.04446915500440518
$ws.Range("K19").Value = 0.4473609171378428
$ws.Range("N19").Value = 1.276474819749449
$ws.Range("O19").Value = 2.284286341898252

$ws.Range("B20").Value = 0.4434815065283146
$ws.Range("C20").Value = 0.06256246358590545
$ws.Range("E20").Value = 0.774881613984121
$ws.Range("F20").Value = 2.396446312242944
$ws.Range("G20").Value = 0.002423709110729057
$ws.Range("I20").Value = 0.4866180566031453
$ws.Range("J20").Value = 0.04436117012868479
$ws.Range("K20").Value = 0.4634556022020035
$ws.Range("N20").Value = 1.272086034570478
$ws.Range("O20").Value = 2.279136623333031

$ws.Range("B21").Value = 0.4931808370883743
$ws.Range("C21").Value = 0.07021052478350498
$ws.Range("E21").Value = 0.8170022093123634
$ws.Range("F21").Value = 2.449626101076092
$ws.Range("G21").Value = 0.002421314298808236
$ws.Range("I21").Value = 0.4807663424773416
$ws.Range("J21").Value = 0.04403021715989297
$ws.Range("K21").Value = 0.5174986689314096
$ws.Range("N21").Value = 1.257883642141891
$ws.Range("O21").Value = 2.263456154342691

$ws.Range("B22").Value = 0.5256533818001117
$ws.Range("C22").Value = 0.07519722248797223
$ws.Range("E22").Value = 0.8448370049441678
$ws.Range("F22").Value = 2.485586145855564
$ws.Range("G22").Value = 0.002419807268647326
$ws.Range("I22").Value = 0.4771862584458972
$ws.Range("J22").Value = 0.04383771338738285
$ws.Range("K22").Value = 0.5527914532360114
$ws.Range("N22").Value = 1.24900654044113
$ws.Range("O22").Value = 2.254418680847522

$ws.Range("B23").Value = 0.5083231763297817
$ws.Range("C23").Value = 0.07253682056455091
$ws.Range("E23").Value = 0.8299529635148417
$ws.Range("F23").Value = 2.466283278670716
$ws.Range("G23").Value = 0.002420606296268576
$ws.Range("I23").Value = 0.4790746368696652
$ws.Range("J23").Value = 0.04393827621639446
$ws.Range("K23").Value = 0.5339577363973547
$ws.Range("N23").Value = 1.253707357977312
$ws.Range("O23").Value = 2.259131262410108

$ws.Range("B24").Value = 0.4426690907295949
$ws.Range("C24").Value = 0.06243726541107719
$ws.Range("E24").Value = 0.7741983415410942
$ws.Range("F24").Value = 2.395597408316718
$ws.Range("G24").Value = 0.002423749237531019
$ws.Range("I24").Value = 0.4867178201046229
$ws.Range("J24").Value = 0.04436698039639708
$ws.Range("K24").Value = 0.4625718750446026
$ws.Range("N24").Value = 1.272325001745212
$ws.Range("O24").Value = 2.2794133243017

$ws.Range("B25").Value = 0.3718318972936174
$ws.Range("C25").Value = 0.05149449886972945
$ws.Range("E25").Value = 0.7153701031286772
$ws.Range("F25").Value = 2.324504934259735
$ws.Range("G25").Value = 0.002427390488150344
$ws.Range("I25").Value = 0.4960091310416672
$ws.Range("J25").Value = 0.04493125997069569
$ws.Range("K25").Value = 0.3854717414335482
$ws.Range("N25").Value = 1.294145741608091
$ws.Range("O25").Value = 2.306476724107114
